$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.979.23'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '3.522.50'
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.35%  '
$ws.Range("D7").Value = '3.520.62'
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("E11").Value = '  -4.50%  '
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").Value = '4.119.60'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("E14").Value = '  -8.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.77%  '
$ws.Range("D16").Value = '3.528.05'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = '65.876.23'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("E19").Value = '  -4.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '420.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.02%  '
$ws.Range("E23").Value = '  -4.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '76.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.81%  '
$ws.Range("D25").Value = '3.664.02'
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("E27").Value = '  -7.09%  '
$ws.Range("E28").Value = '  -2.53%  '
$ws.Range("E29").Value = '  -6.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '3.529.34'
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.154'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.53%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  -8.06%  '
$ws.Range("E37").Value = '  -4.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '177.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("E39").Value = '  -5.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.73%  '
$ws.Range("E41").Value = '  -5.13%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.858'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.53%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.96'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("E45").Value = '  -8.37%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("E49").Value = '  -2.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.905'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.15%  '
